$d = $word.ActiveDocument

$d.Content.Find.Execute("50+5=", $true, $false, $false, $false, $false, $true, 1, $false, "19+39=", 2) | Out-Null
$d.Content.Find.Execute("63+8=", $true, $false, $false, $false, $false, $true, 1, $false, "84-2=", 2) | Out-Null
$d.Content.Find.Execute("55-37=", $true, $false, $false, $false, $false, $true, 1, $false, "36-12=", 2) | Out-Null
$d.Content.Find.Execute("90-12=", $true, $false, $false, $false, $false, $true, 1, $false, "92-18=", 2) | Out-Null
$d.Content.Find.Execute("78-24=", $true, $false, $false, $false, $false, $true, 1, $false, "11+55=", 2) | Out-Null
$d.Content.Find.Execute("98-57=", $true, $false, $false, $false, $false, $true, 1, $false, "33-28=", 2) | Out-Null
$d.Content.Find.Execute("17+1=", $true, $false, $false, $false, $false, $true, 1, $false, "72-12=", 2) | Out-Null
$d.Content.Find.Execute("79-35=", $true, $false, $false, $false, $false, $true, 1, $false, "88-17=", 2) | Out-Null
$d.Content.Find.Execute("73+16=", $true, $false, $false, $false, $false, $true, 1, $false, "1+79=", 2) | Out-Null
$d.Content.Find.Execute("33+53=", $true, $false, $false, $false, $false, $true, 1, $false, "27+3=", 2) | Out-Null
$d.Content.Find.Execute("79-28=", $true, $false, $false, $false, $false, $true, 1, $false, "79+18=", 2) | Out-Null
$d.Content.Find.Execute("48-12=", $true, $false, $false, $false, $false, $true, 1, $false, "82-66=", 2) | Out-Null
$d.Content.Find.Execute("47+10=", $true, $false, $false, $false, $false, $true, 1, $false, "79-42=", 2) | Out-Null
$d.Content.Find.Execute("66-46=", $true, $false, $false, $false, $false, $true, 1, $false, "23+76=", 2) | Out-Null
$d.Content.Find.Execute("97-76=", $true, $false, $false, $false, $false, $true, 1, $false, "41-27=", 2) | Out-Null
$d.Content.Find.Execute("86+4=", $true, $false, $false, $false, $false, $true, 1, $false, "9+76=", 2) | Out-Null
$d.Content.Find.Execute("26-5=", $true, $false, $false, $false, $false, $true, 1, $false, "2+46=", 2) | Out-Null
$d.Content.Find.Execute("92-27=", $true, $false, $false, $false, $false, $true, 1, $false, "87-6=", 2) | Out-Null
$d.Content.Find.Execute("41+3=", $true, $false, $false, $false, $false, $true, 1, $false, "93-2=", 2) | Out-Null
$d.Content.Find.Execute("69+0=", $true, $false, $false, $false, $false, $true, 1, $false, "61-58=", 2) | Out-Null
$d.Content.Find.Execute("96-74=", $true, $false, $false, $false, $false, $true, 1, $false, "73+12=", 2) | Out-Null
$d.Content.Find.Execute("49+14=", $true, $false, $false, $false, $false, $true, 1, $false, "42+43=", 2) | Out-Null
$d.Content.Find.Execute("74-68=", $true, $false, $false, $false, $false, $true, 1, $false, "42+51=", 2) | Out-Null
$d.Content.Find.Execute("24-21=", $true, $false, $false, $false, $false, $true, 1, $false, "87-22=", 2) | Out-Null
$d.Content.Find.Execute("60-19=", $true, $false, $false, $false, $false, $true, 1, $false, "63-42=", 2) | Out-Null
$d.Content.Find.Execute("35+48=", $true, $false, $false, $false, $false, $true, 1, $false, "65-35=", 2) | Out-Null
$d.Content.Find.Execute("36+17=", $true, $false, $false, $false, $false, $true, 1, $false, "64+15=", 2) | Out-Null
$d.Content.Find.Execute("88-56=", $true, $false, $false, $false, $false, $true, 1, $false, "61-44=", 2) | Out-Null
$d.Content.Find.Execute("90-7=", $true, $false, $false, $false, $false, $true, 1, $false, "39+53=", 2) | Out-Null
$d.Content.Find.Execute("76+4=", $true, $false, $false, $false, $false, $true, 1, $false, "80+1=", 2) | Out-Null
$d.Content.Find.Execute("65+34=", $true, $false, $false, $false, $false, $true, 1, $false, "18+45=", 2) | Out-Null
$d.Content.Find.Execute("9+72=", $true, $false, $false, $false, $false, $true, 1, $false, "93-38=", 2) | Out-Null
$d.Content.Find.Execute("52-20=", $true, $false, $false, $false, $false, $true, 1, $false, "80+6=", 2) | Out-Null
$d.Content.Find.Execute("52+24=", $true, $false, $false, $false, $false, $true, 1, $false, "96-40=", 2) | Out-Null
$d.Content.Find.Execute("56+34=", $true, $false, $false, $false, $false, $true, 1, $false, "49-22=", 2) | Out-Null
$d.Content.Find.Execute("89-86=", $true, $false, $false, $false, $false, $true, 1, $false, "80+16=", 2) | Out-Null
$d.Content.Find.Execute("47+27=", $true, $false, $false, $false, $false, $true, 1, $false, "89-88=", 2) | Out-Null
$d.Content.Find.Execute("41-4=", $true, $false, $false, $false, $false, $true, 1, $false, "33+34=", 2) | Out-Null
$d.Content.Find.Execute("72+3=", $true, $false, $false, $false, $false, $true, 1, $false, "75-72=", 2) | Out-Null
$d.Content.Find.Execute("19+27=", $true, $false, $false, $false, $false, $true, 1, $false, "63-2=", 2) | Out-Null
$d.Content.Find.Execute("58-25=", $true, $false, $false, $false, $false, $true, 1, $false, "99-56=", 2) | Out-Null
$d.Content.Find.Execute("74-5=", $true, $false, $false, $false, $false, $true, 1, $false, "74+25=", 2) | Out-Null
$d.Content.Find.Execute("72-65=", $true, $false, $false, $false, $false, $true, 1, $false, "50-41=", 2) | Out-Null
$d.Content.Find.Execute("71-3=", $true, $false, $false, $false, $false, $true, 1, $false, "94+5=", 2) | Out-Null
$d.Content.Find.Execute("43-36=", $true, $false, $false, $false, $false, $true, 1, $false, "64-29=", 2) | Out-Null
$d.Content.Find.Execute("14-6=", $true, $false, $false, $false, $false, $true, 1, $false, "25-10=", 2) | Out-Null
$d.Content.Find.Execute("68+24=", $true, $false, $false, $false, $false, $true, 1, $false, "76-27=", 2) | Out-Null
$d.Content.Find.Execute("54-46=", $true, $false, $false, $false, $false, $true, 1, $false, "6+66=", 2) | Out-Null
$d.Content.Find.Execute("22-0=", $true, $false, $false, $false, $false, $true, 1, $false, "41+7=", 2) | Out-Null
$d.Content.Find.Execute("8+28=", $true, $false, $false, $false, $false, $true, 1, $false, "69+15=", 2) | Out-Null
$d.Content.Find.Execute("83-20=", $true, $false, $false, $false, $false, $true, 1, $false, "92-40=", 2) | Out-Null
$d.Content.Find.Execute("18+52=", $true, $false, $false, $false, $false, $true, 1, $false, "90+4=", 2) | Out-Null
$d.Content.Find.Execute("74-25=", $true, $false, $false, $false, $false, $true, 1, $false, "43-35=", 2) | Out-Null
$d.Content.Find.Execute("52+4=", $true, $false, $false, $false, $false, $true, 1, $false, "44+34=", 2) | Out-Null
$d.Content.Find.Execute("82-57=", $true, $false, $false, $false, $false, $true, 1, $false, "96-86=", 2) | Out-Null
$d.Content.Find.Execute("55-31=", $true, $false, $false, $false, $false, $true, 1, $false, "80-18=", 2) | Out-Null
$d.Content.Find.Execute("26+50=", $true, $false, $false, $false, $false, $true, 1, $false, "87-19=", 2) | Out-Null
$d.Content.Find.Execute("22+75=", $true, $false, $false, $false, $false, $true, 1, $false, "35+6=", 2) | Out-Null
$d.Content.Find.Execute("34+17=", $true, $false, $false, $false, $false, $true, 1, $false, "9+53=", 2) | Out-Null
$d.Content.Find.Execute("93-22=", $true, $false, $false, $false, $false, $true, 1, $false, "23+69=", 2) | Out-Null
$d.Content.Find.Execute("18+79=", $true, $false, $false, $false, $false, $true, 1, $false, "28+27=", 2) | Out-Null
$d.Content.Find.Execute("73-37=", $true, $false, $false, $false, $false, $true, 1, $false, "87-58=", 2) | Out-Null
$d.Content.Find.Execute("88-33=", $true, $false, $false, $false, $false, $true, 1, $false, "41-34=", 2) | Out-Null
$d.Content.Find.Execute("19+47=", $true, $false, $false, $false, $false, $true, 1, $false, "50-23=", 2) | Out-Null
$d.Content.Find.Execute("2+36=", $true, $false, $false, $false, $false, $true, 1, $false, "23-3=", 2) | Out-Null
$d.Content.Find.Execute("30+35=", $true, $false, $false, $false, $false, $true, 1, $false, "42-6=", 2) | Out-Null
$d.Content.Find.Execute("85-0=", $true, $false, $false, $false, $false, $true, 1, $false, "25+26=", 2) | Out-Null
$d.Content.Find.Execute("50+37=", $true, $false, $false, $false, $false, $true, 1, $false, "28-4=", 2) | Out-Null
$d.Content.Find.Execute("76-13=", $true, $false, $false, $false, $false, $true, 1, $false, "82-18=", 2) | Out-Null
$d.Content.Find.Execute("60-1=", $true, $false, $false, $false, $false, $true, 1, $false, "49+12=", 2) | Out-Null
$d.Content.Find.Execute("67-32=", $true, $false, $false, $false, $false, $true, 1, $false, "29-14=", 2) | Out-Null
$d.Content.Find.Execute("1+80=", $true, $false, $false, $false, $false, $true, 1, $false, "72-28=", 2) | Out-Null
$d.Content.Find.Execute("57+8=", $true, $false, $false, $false, $false, $true, 1, $false, "70-33=", 2) | Out-Null
$d.Content.Find.Execute("60+3=", $true, $false, $false, $false, $false, $true, 1, $false, "37+50=", 2) | Out-Null
$d.Content.Find.Execute("63-7=", $true, $false, $false, $false, $false, $true, 1, $false, "46+3=", 2) | Out-Null
$d.Content.Find.Execute("95+3=", $true, $false, $false, $false, $false, $true, 1, $false, "37+58=", 2) | Out-Null
$d.Content.Find.Execute("55-48=", $true, $false, $false, $false, $false, $true, 1, $false, "16+24=", 2) | Out-Null
$d.Content.Find.Execute("82+4=", $true, $false, $false, $false, $false, $true, 1, $false, "74+24=", 2) | Out-Null
$d.Content.Find.Execute("18+53=", $true, $false, $false, $false, $false, $true, 1, $false, "67+27=", 2) | Out-Null
$d.Content.Find.Execute("35-6=", $true, $false, $false, $false, $false, $true, 1, $false, "6+47=", 2) | Out-Null
$d.Content.Find.Execute("54+8=", $true, $false, $false, $false, $false, $true, 1, $false, "59-5=", 2) | Out-Null
$d.Content.Find.Execute("30+19=", $true, $false, $false, $false, $false, $true, 1, $false, "86-34=", 2) | Out-Null
$d.Content.Find.Execute("29+8=", $true, $false, $false, $false, $false, $true, 1, $false, "23-9=", 2) | Out-Null
$d.Content.Find.Execute("16+15=", $true, $false, $false, $false, $false, $true, 1, $false, "94-71=", 2) | Out-Null
$d.Content.Find.Execute("53-25=", $true, $false, $false, $false, $false, $true, 1, $false, "32-20=", 2) | Out-Null
$d.Content.Find.Execute("32-24=", $true, $false, $false, $false, $false, $true, 1, $false, "78+11=", 2) | Out-Null
$d.Content.Find.Execute("65+11=", $true, $false, $false, $false, $false, $true, 1, $false, "5+40=", 2) | Out-Null
$d.Content.Find.Execute("64-22=", $true, $false, $false, $false, $false, $true, 1, $false, "21+58=", 2) | Out-Null
$d.Content.Find.Execute("92-34=", $true, $false, $false, $false, $false, $true, 1, $false, "8-0=", 2) | Out-Null
$d.Content.Find.Execute("41+49=", $true, $false, $false, $false, $false, $true, 1, $false, "51-25=", 2) | Out-Null
$d.Content.Find.Execute("50+20=", $true, $false, $false, $false, $false, $true, 1, $false, "15+31=", 2) | Out-Null
$d.Content.Find.Execute("43-7=", $true, $false, $false, $false, $false, $true, 1, $false, "45-20=", 2) | Out-Null
$d.Content.Find.Execute("90-78=", $true, $false, $false, $false, $false, $true, 1, $false, "83-58=", 2) | Out-Null
$d.Content.Find.Execute("98-76=", $true, $false, $false, $false, $false, $true, 1, $false, "6+76=", 2) | Out-Null
$d.Content.Find.Execute("92-6=", $true, $false, $false, $false, $false, $true, 1, $false, "55+44=", 2) | Out-Null
$d.Content.Find.Execute("94-52=", $true, $false, $false, $false, $false, $true, 1, $false, "54-38=", 2) | Out-Null
$d.Content.Find.Execute("96-12=", $true, $false, $false, $false, $false, $true, 1, $false, "71+1=", 2) | Out-Null
$d.Content.Find.Execute("93+3=", $true, $false, $false, $false, $false, $true, 1, $false, "63-33=", 2) | Out-Null
$d.Content.Find.Execute("6+65=", $true, $false, $false, $false, $false, $true, 1, $false, "42-29=", 2) | Out-Null
$d.Content.Find.Execute("4-3=", $true, $false, $false, $false, $false, $true, 1, $false, "94-64=", 2) | Out-Null
